$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.167466163635254
$ws.Range("B1").Value = 2.491630554199219
$ws.Range("C1").Value = 6.674103260040283
$ws.Range("D1").Value = 2.05537748336792
$ws.Range("E1").Value = 1.20969033241272
